$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = "model_26_4_24"
$ws.Range("B2").Value2 = 0.9999770883131851
$ws.Range("C2").Value2 = 0.999055431844164
$ws.Range("D2").Value2 = 0.9998938733699413
$ws.Range("E2").Value2 = 0.9998502534770579
$ws.Range("F2").Value2 = 0.9999358483145446
$ws.Range("G2").Value2 = 0.00002138705393491294
$ws.Range("H2").Value2 = 0.0008817129117227806
$ws.Range("I2").Value2 = 0.00002533296306992014
$ws.Range("J2").Value2 = 0.0001079511249075915
$ws.Range("K2").Value2 = 0.00006664204398875581
$ws.Range("L2").Value2 = 0.0002291822534573157
$ws.Range("M2").Value2 = 0.004624613922795387
$ws.Range("N2").Value2 = 1.000026184784931
$ws.Range("O2").Value2 = 0.004821493311759351
$ws.Range("P2").Value2 = 111.5054495506177
$ws.Range("Q2").Value2 = 166.3548616696867

$ws.Range("A3").Value2 = "model_26_4_23"
$ws.Range("B3").Value2 = 0.9999774029545943
$ws.Range("C3").Value2 = 0.9990540245338234
$ws.Range("D3").Value2 = 0.9998952704433834
$ws.Range("E3").Value2 = 0.9998524651823261
$ws.Range("F3").Value2 = 0.9999367774092849
$ws.Range("G3").Value2 = 0.00002109334999065019
$ws.Range("H3").Value2 = 0.0008830265741519294
$ws.Range("I3").Value2 = 0.00002499947457701045
$ws.Range("J3").Value2 = 0.0001063567234686225
$ws.Range("K3").Value2 = 0.00006567688193402271
$ws.Range("L3").Value2 = 0.0002270011686536036
$ws.Range("M3").Value2 = 0.004592749720009811
$ws.Range("N3").Value2 = 1.000025825194749
$ws.Range("O3").Value2 = 0.004788272583893199
$ws.Range("P3").Value2 = 111.533105466955
$ws.Range("Q3").Value2 = 166.382517586024

$ws.Range("A4").Value2 = "model_26_4_22"
$ws.Range("B4").Value2 = 0.9999777396621201
$ws.Range("C4").Value2 = 0.9990524471718454
$ws.Range("D4").Value2 = 0.9998966037025505
$ws.Range("E4").Value2 = 0.9998548936539946
$ws.Range("F4").Value2 = 0.999937776179855
$ws.Range("G4").Value2 = 0.00002077904829496498
$ws.Range("H4").Value2 = 0.0008844989723202503
$ws.Range("I4").Value2 = 0.00002468121887415554
$ws.Range("J4").Value2 = 0.0001046060567868858
$ws.Range("K4").Value2 = 0.0000646393392444401
$ws.Range("L4").Value2 = 0.0002245398291589179
$ws.Range("M4").Value2 = 0.004558404139056231
$ws.Range("N4").Value2 = 1.000025440386148
$ws.Range("O4").Value2 = 0.004752464840453272
$ws.Range("P4").Value2 = 111.5631307448587
$ws.Range("Q4").Value2 = 166.4125428639277

$ws.Range("A5").Value2 = "model_26_4_21"
$ws.Range("B5").Value2 = 0.9999780965178211
$ws.Range("C5").Value2 = 0.9990506767191827
$ws.Range("D5").Value2 = 0.9998978764695279
$ws.Range("E5").Value2 = 0.9998575548795027
$ws.Range("F5").Value2 = 0.999938841656172
$ws.Range("G5").Value2 = 0.00002044593916226536
$ws.Range("H5").Value2 = 0.0008861516121669859
$ws.Range("I5").Value2 = 0.00002437740296275882
$ws.Range("J5").Value2 = 0.0001026875996395382
$ws.Range("K5").Value2 = 0.00006353250130114849
$ws.Range("L5").Value2 = 0.000221810454201872
$ws.Range("M5").Value2 = 0.004521718607152081
$ws.Range("N5").Value2 = 1.000025032551062
$ws.Range("O5").Value2 = 0.004714217529506443
$ws.Range("P5").Value2 = 111.5954525383325
$ws.Range("Q5").Value2 = 166.4448646574015

$ws.Range("A6").Value2 = "model_26_4_20"
$ws.Range("B6").Value2 = 0.9999784723780656
$ws.Range("C6").Value2 = 0.9990486857344907
$ws.Range("D6").Value2 = 0.9998989698760368
$ws.Range("E6").Value2 = 0.999860468791306
$ws.Range("F6").Value2 = 0.9999399783395145
$ws.Range("G6").Value2 = 0.0000200950901224114
$ws.Range("H6").Value2 = 0.0008880101089828233
$ws.Range("I6").Value2 = 0.00002411640130185745
$ws.Range("J6").Value2 = 0.0001005869828714401
$ws.Range("K6").Value2 = 0.00006235169208664879
$ws.Range("L6").Value2 = 0.0002187596468423212
$ws.Range("M6").Value2 = 0.004482754747073656
$ws.Range("N6").Value2 = 1.000024602996497
$ws.Range("O6").Value2 = 0.004673594897238169
$ws.Range("P6").Value2 = 111.630070090511
$ws.Range("Q6").Value2 = 166.47948220958

$ws.Range("A7").Value2 = "model_26_4_19"
$ws.Range("B7").Value2 = 0.999978864846521
$ws.Range("C7").Value2 = 0.9990464461045662
$ws.Range("D7").Value2 = 0.9998998113882895
$ws.Range("E7").Value2 = 0.9998636528209418
$ws.Range("F7").Value2 = 0.9999411775350597
$ws.Range("G7").Value2 = 0.00001972873804668808
$ws.Range("H7").Value2 = 0.0008901007052090147
$ws.Range("I7").Value2 = 0.00002391552807324117
$ws.Range("J7").Value2 = 0.00009829164022060452
$ws.Range("K7").Value2 = 0.00006110594395535005
$ws.Range("L7").Value2 = 0.0002175892725284323
$ws.Range("M7").Value2 = 0.004441704407847069
$ws.Range("N7").Value2 = 1.000024154461119
$ws.Range("O7").Value2 = 0.004630796960084789
$ws.Range("P7").Value2 = 111.6668684022327
$ws.Range("Q7").Value2 = 166.5162805213017

$ws.Range("A8").Value2 = "model_26_4_18"
$ws.Range("B8").Value2 = 0.9999792700978076
$ws.Range("C8").Value2 = 0.999043918867406
$ws.Range("D8").Value2 = 0.9999003246950727
$ws.Range("E8").Value2 = 0.999867117529697
$ws.Range("F8").Value2 = 0.9999424408889765
$ws.Range("G8").Value2 = 0.0000193504537590877
$ws.Range("H8").Value2 = 0.0008924597701651207
$ws.Range("I8").Value2 = 0.00002379299914930285
$ws.Range("J8").Value2 = 0.00009579395813582842
$ws.Range("K8").Value2 = 0.00005979354683428757
$ws.Range("L8").Value2 = 0.0002191593212289637
$ws.Range("M8").Value2 = 0.00439891506613707
$ws.Range("N8").Value2 = 1.000023691316791
$ws.Range("O8").Value2 = 0.004586185987512047
$ws.Range("P8").Value2 = 111.7055893773628
$ws.Range("Q8").Value2 = 166.5550014964318

$ws.Range("A9").Value2 = "model_26_4_17"
$ws.Range("B9").Value2 = 0.9999796799503229
$ws.Range("C9").Value2 = 0.9990410679854588
$ws.Range("D9").Value2 = 0.9999003026219062
$ws.Range("E9").Value2 = 0.9998708831131409
$ws.Range("F9").Value2 = 0.9999437432460939
$ws.Range("G9").Value2 = 0.0000189678744264612
$ws.Range("H9").Value2 = 0.0008951209433236155
$ws.Range("I9").Value2 = 0.00002379826812573634
$ws.Range("J9").Value2 = 0.00009307937778561767
$ws.Range("K9").Value2 = 0.00005844063241455244
$ws.Range("L9").Value2 = 0.0002230054978268744
$ws.Range("M9").Value2 = 0.004355212328516394
$ws.Range("N9").Value2 = 1.000023222913917
$ws.Range("O9").Value2 = 0.004540622733873722
$ws.Range("P9").Value2 = 111.7455276590133
$ws.Range("Q9").Value2 = 166.5949397780823

$ws.Range("A10").Value2 = "model_26_4_16"
$ws.Range("B10").Value2 = 0.9999800883126342
$ws.Range("C10").Value2 = 0.9990378467356966
$ws.Range("D10").Value2 = 0.9998995803196974
$ws.Range("E10").Value2 = 0.9998749607144344
$ws.Range("F10").Value2 = 0.999945076640159
$ws.Range("G10").Value2 = 0.00001858668613389017
$ws.Range("H10").Value2 = 0.00089812783858001
$ws.Range("I10").Value2 = 0.00002397068531424595
$ws.Range("J10").Value2 = 0.00009013986615015843
$ws.Range("K10").Value2 = 0.0000570554761975713
$ws.Range("L10").Value2 = 0.0002272567156107493
$ws.Range("M10").Value2 = 0.00431122791486256
$ws.Range("N10").Value2 = 1.000022756214132
$ws.Range("O10").Value2 = 0.004494765812670354
$ws.Range("P10").Value2 = 111.7861300656099
$ws.Range("Q10").Value2 = 166.6355421846789

$ws.Range("A11").Value2 = "model_26_4_15"
$ws.Range("B11").Value2 = 0.999980482414467
$ws.Range("C11").Value2 = 0.9990342025616544
$ws.Range("D11").Value2 = 0.9998979474680677
$ws.Range("E11").Value2 = 0.9998793574188437
$ws.Range("F11").Value2 = 0.9999464130762852
$ws.Range("G11").Value2 = 0.00001821880937201312
$ws.Range("H11").Value2 = 0.0009015295150876449
$ws.Range("I11").Value2 = 0.0000243604552524
$ws.Range("J11").Value2 = 0.00008697031551524024
$ws.Range("K11").Value2 = 0.00005566715982714954
$ws.Range("L11").Value2 = 0.0002319272688400162
$ws.Range("M11").Value2 = 0.00426834972466094
$ws.Range("N11").Value2 = 1.000022305812038
$ws.Range("O11").Value2 = 0.004450062209141761
$ws.Range("P11").Value2 = 111.8261120311552
$ws.Range("Q11").Value2 = 166.6755241502243

$ws.Range("A12").Value2 = "model_26_4_14"
$ws.Range("B12").Value2 = 0.9999808465910265
$ws.Range("C12").Value2 = 0.9990300731050227
$ws.Range("D12").Value2 = 0.9998950680003733
$ws.Range("E12").Value2 = 0.9998840705236491
$ws.Range("F12").Value2 = 0.9999477192925015
$ws.Range("G12").Value2 = 0.00001787886653926069
$ws.Range("H12").Value2 = 0.0009053841815911783
$ws.Range("I12").Value2 = 0.00002504779874687767
$ws.Range("J12").Value2 = 0.00008357267425075034
$ws.Range("K12").Value2 = 0.000054310236498814
$ws.Range("L12").Value2 = 0.0002370921123333001
$ws.Range("M12").Value2 = 0.004228340873115683
$ws.Range("N12").Value2 = 1.000021889610256
$ws.Range("O12").Value2 = 0.004408350097956482
$ws.Range("P12").Value2 = 111.8637823657742
$ws.Range("Q12").Value2 = 166.7131944848432

$ws.Range("A13").Value2 = "model_26_4_13"
$ws.Range("B13").Value2 = 0.9999811606632489
$ws.Range("C13").Value2 = 0.9990253827117181
$ws.Range("D13").Value2 = 0.9998905872706882
$ws.Range("E13").Value2 = 0.9998890997140671
$ws.Range("F13").Value2 = 0.9999489493048983
$ws.Range("G13").Value2 = 0.00001758569390580093
$ws.Range("H13").Value2 = 0.0009097624578564802
$ws.Range("I13").Value2 = 0.00002611737157300519
$ws.Range("J13").Value2 = 0.00007994716928186341
$ws.Range("K13").Value2 = 0.0000530324752104613
$ws.Range("L13").Value2 = 0.0002427793396746841
$ws.Range("M13").Value2 = 0.004193530005353596
$ws.Range("N13").Value2 = 1.000021530670573
$ws.Range("O13").Value2 = 0.004372057259485334
$ws.Range("P13").Value2 = 111.8968496654403
$ws.Range("Q13").Value2 = 166.7462617845093

$ws.Range("A14").Value2 = "model_26_4_12"
$ws.Range("B14").Value2 = 0.9999813958546965
$ws.Range("C14").Value2 = 0.9990200508365611
$ws.Range("D14").Value2 = 0.999883997871661
$ws.Range("E14").Value2 = 0.9998944216961638
$ws.Range("F14").Value2 = 0.9999500390332916
$ws.Range("G14").Value2 = 0.000017366153013198
$ws.Range("H14").Value2 = 0.0009147395292733344
$ws.Range("I14").Value2 = 0.00002769029443049799
$ws.Range("J14").Value2 = 0.00007611059302763794
$ws.Range("K14").Value2 = 0.00005190044372906797
$ws.Range("L14").Value2 = 0.0002490275001991013
$ws.Range("M14").Value2 = 0.004167271650996369
$ws.Range("N14").Value2 = 1.000021261880347
$ws.Range("O14").Value2 = 0.004344681032620832
$ws.Range("P14").Value2 = 111.9219749505713
$ws.Range("Q14").Value2 = 166.7713870696403

$ws.Range("A15").Value2 = "model_26_4_11"
$ws.Range("B15").Value2 = 0.9999815215923202
$ws.Range("C15").Value2 = 0.999013976992698
$ws.Range("D15").Value2 = 0.999874776944123
$ws.Range("E15").Value2 = 0.9999000102440502
$ws.Range("F15").Value2 = 0.9999509187049104
$ws.Range("G15").Value2 = 0.00001724878246073583
$ws.Range("H15").Value2 = 0.0009204091959086695
$ws.Range("I15").Value2 = 0.00002989137644604828
$ws.Range("J15").Value2 = 0.00007208185153118515
$ws.Range("K15").Value2 = 0.00005098662339369545
$ws.Range("L15").Value2 = 0.0002558829700206974
$ws.Range("M15").Value2 = 0.004153165354369583
$ws.Range("N15").Value2 = 1.000021118180205
$ws.Range("O15").Value2 = 0.004329974201742585
$ws.Range("P15").Value2 = 111.9355379979285
$ws.Range("Q15").Value2 = 166.7849501169975

$ws.Range("A16").Value2 = "model_26_4_10"
$ws.Range("B16").Value2 = 0.9999814908832334
$ws.Range("C16").Value2 = 0.9990070341228666
$ws.Range("D16").Value2 = 0.9998621062973351
$ws.Range("E16").Value2 = 0.9999058293091049
$ws.Range("F16").Value2 = 0.9999514801816524
$ws.Range("G16").Value2 = 0.00001727744804533376
$ws.Range("H16").Value2 = 0.0009268900601394649
$ws.Range("I16").Value2 = 0.00003291592388502341
$ws.Range("J16").Value2 = 0.00006788693196824825
$ws.Range("K16").Value2 = 0.00005040335021126808
$ws.Range("L16").Value2 = 0.0002633996439785122
$ws.Range("M16").Value2 = 0.004156614974391273
$ws.Range("N16").Value2 = 1.000021153276305
$ws.Range("O16").Value2 = 0.004333570679230295
$ws.Range("P16").Value2 = 111.9322169761124
$ws.Range("Q16").Value2 = 166.7816290951814

$ws.Range("A17").Value2 = "model_26_4_9"
$ws.Range("B17").Value2 = 0.9999812433614558
$ws.Range("C17").Value2 = 0.9989990938122534
$ws.Range("D17").Value2 = 0.9998450565137121
$ws.Range("E17").Value2 = 0.9999117697557993
$ws.Range("F17").Value2 = 0.9999515843371958
$ws.Range("G17").Value2 = 0.00001750849875976973
$ws.Range("H17").Value2 = 0.0009343019915575142
$ws.Range("I17").Value2 = 0.00003698579342324472
$ws.Range("J17").Value2 = 0.00006360450930820754
$ws.Range("K17").Value2 = 0.00005029515136572613
$ws.Range("L17").Value2 = 0.0002716320352672009
$ws.Range("M17").Value2 = 0.004184315805453711
$ws.Range("N17").Value2 = 1.000021436158336
$ws.Range("O17").Value2 = 0.004362450792019686
$ws.Range("P17").Value2 = 111.905648303012
$ws.Range("Q17").Value2 = 166.755060422081

$ws.Range("A18").Value2 = "model_26_4_8"
$ws.Range("B18").Value2 = 0.9999807059078945
$ws.Range("C18").Value2 = 0.9989899884969244
$ws.Range("D18").Value2 = 0.9998224813545205
$ws.Range("E18").Value2 = 0.9999177180857919
$ws.Range("F18").Value2 = 0.9999510567993097
$ws.Range("G18").Value2 = 0.00001801018806773756
$ws.Range("H18").Value2 = 0.0009428014037399272
$ws.Range("I18").Value2 = 0.00004237459804073175
$ws.Range("J18").Value2 = 0.00005931640363862385
$ws.Range("K18").Value2 = 0.00005084316819113586
$ws.Range("L18").Value2 = 0.0002806456077365049
$ws.Range("M18").Value2 = 0.004243841192568068
$ws.Range("N18").Value2 = 1.000022050390978
$ws.Range("O18").Value2 = 0.00442451029809804
$ws.Range("P18").Value2 = 111.8491459128493
$ws.Range("Q18").Value2 = 166.6985580319183

$ws.Range("A19").Value2 = "model_26_4_7"
$ws.Range("B19").Value2 = 0.9999797868245517
$ws.Range("C19").Value2 = 0.9989795279245721
$ws.Range("D19").Value2 = 0.9997929472032119
$ws.Range("E19").Value2 = 0.9999235317198009
$ws.Range("F19").Value2 = 0.9999496800908008
$ws.Range("G19").Value2 = 0.00001886811202518873
$ws.Range("H19").Value2 = 0.0009525658888648454
$ws.Range("I19").Value2 = 0.00004942454925456331
$ws.Range("J19").Value2 = 0.00005512539927510625
$ws.Range("K19").Value2 = 0.00005227332031206835
$ws.Range("L19").Value2 = 0.0002904650161404944
$ws.Range("M19").Value2 = 0.004343744010089536
$ws.Range("N19").Value2 = 1.000023100771941
$ws.Range("O19").Value2 = 0.00452866618538874
$ws.Range("P19").Value2 = 111.7560745104973
$ws.Range("Q19").Value2 = 166.6054866295663

$ws.Range("A20").Value2 = "model_26_4_6"
$ws.Range("B20").Value2 = 0.9999783656821398
$ws.Range("C20").Value2 = 0.9989674715562625
$ws.Range("D20").Value2 = 0.9997546479161296
$ws.Range("E20").Value2 = 0.9999289848279875
$ws.Range("F20").Value2 = 0.9999471710041884
$ws.Range("G20").Value2 = 0.00002019468608563559
$ws.Range("H20").Value2 = 0.0009638199794684228
$ws.Range("I20").Value2 = 0.00005856678268574159
$ws.Range("J20").Value2 = 0.0000511942952239276
$ws.Range("K20").Value2 = 0.00005487980927967783
$ws.Range("L20").Value2 = 0.0003011561218065388
$ws.Range("M20").Value2 = 0.004493849806750955
$ws.Range("N20").Value2 = 1.000024724934697
$ws.Range("O20").Value2 = 0.004685162296575869
$ws.Range("P20").Value2 = 111.6201821064647
$ws.Range("Q20").Value2 = 166.4695942255337

$ws.Range("A21").Value2 = "model_26_4_5"
$ws.Range("B21").Value2 = 0.999976291709006
$ws.Range("C21").Value2 = 0.998953546348885
$ws.Range("D21").Value2 = 0.9997053133861391
$ws.Range("E21").Value2 = 0.9999337812685649
$ws.Range("F21").Value2 = 0.9999431663908306
$ws.Range("G21").Value2 = 0.00002213064896919195
$ws.Range("H21").Value2 = 0.0009768185493094113
$ws.Range("I21").Value2 = 0.00007034318438274349
$ws.Range("J21").Value2 = 0.00004773657783785593
$ws.Range("K21").Value2 = 0.00005903988111029971
$ws.Range("L21").Value2 = 0.0003127656230357724
$ws.Range("M21").Value2 = 0.004704322370883181
$ws.Range("N21").Value2 = 1.000027095189707
$ws.Range("O21").Value2 = 0.004904595113501477
$ws.Range("P21").Value2 = 111.4370941580352
$ws.Range("Q21").Value2 = 166.2865062771043

$ws.Range("A22").Value2 = "model_26_4_4"
$ws.Range("B22").Value2 = 0.9999733744898419
$ws.Range("C22").Value2 = 0.9989374296448014
$ws.Range("D22").Value2 = 0.9996421754497713
$ws.Range("E22").Value2 = 0.9999375130766669
$ws.Range("F22").Value2 = 0.9999372071128827
$ws.Range("G22").Value2 = 0.00002485374500775344
$ws.Range("H22").Value2 = 0.0009918627851298986
$ws.Range("I22").Value2 = 0.00008541452895885279
$ws.Range("J22").Value2 = 0.00004504634587366725
$ws.Range("K22").Value2 = 0.00006523049730884281
$ws.Range("L22").Value2 = 0.0003254521510651548
$ws.Range("M22").Value2 = 0.004985353047453454
$ws.Range("N22").Value2 = 1.000030429154466
$ws.Range("O22").Value2 = 0.005197589847787045
$ws.Range("P22").Value2 = 111.2050042244781
$ws.Range("Q22").Value2 = 166.0544163435472

$ws.Range("A23").Value2 = "model_26_4_3"
$ws.Range("B23").Value2 = 0.9999693754300233
$ws.Range("C23").Value2 = 0.9989187318865926
$ws.Range("D23").Value2 = 0.9995617880591434
$ws.Range("E23").Value2 = 0.9999396670946252
$ws.Range("F23").Value2 = 0.9999287186593157
$ws.Range("G23").Value2 = 0.00002858669180999295
$ws.Range("H23").Value2 = 0.001009316321681059
$ws.Range("I23").Value2 = 0.0001046034054636142
$ws.Range("J23").Value2 = 0.00004349353077580421
$ws.Range("K23").Value2 = 0.00007404847133386944
$ws.Range("L23").Value2 = 0.000339192651461637
$ws.Range("M23").Value2 = 0.005346652392852274
$ws.Range("N23").Value2 = 1.000034999508545
$ws.Range("O23").Value2 = 0.005574270454312243
$ws.Range("P23").Value2 = 110.9251385395259
$ws.Range("Q23").Value2 = 165.7745506585949

$ws.Range("A24").Value2 = "model_26_4_2"
$ws.Range("B24").Value2 = 0.9999639921241794
$ws.Range("C24").Value2 = 0.9988969698982162
$ws.Range("D24").Value2 = 0.9994598106029945
$ws.Range("E24").Value2 = 0.9999395682366251
$ws.Range("F24").Value2 = 0.9999169679362826
$ws.Range("G24").Value2 = 0.00003361177151540036
$ws.Range("H24").Value2 = 0.001029630182589572
$ws.Range("I24").Value2 = 0.0001289459397470272
$ws.Range("J24").Value2 = 0.00004356479675322174
$ws.Range("K24").Value2 = 0.00008625535562245583
$ws.Range("L24").Value2 = 0.0003539585923812019
$ws.Range("M24").Value2 = 0.005797565999227638
$ws.Range("N24").Value2 = 1.000041151858081
$ws.Range("O24").Value2 = 0.00604438038643085
$ws.Range("P24").Value2 = 110.6012684192611
$ws.Range("Q24").Value2 = 165.4506805383301

$ws.Range("A25").Value2 = "model_26_4_1"
$ws.Range("B25").Value2 = 0.9999568478852913
$ws.Range("C25").Value2 = 0.9988715730542866
$ws.Range("D25").Value2 = 0.9993309101389117
$ws.Range("E25").Value2 = 0.9999363309356666
$ws.Range("F25").Value2 = 0.9999010349321749
$ws.Range("G25").Value2 = 0.00004028060492166812
$ws.Range("H25").Value2 = 0.001053337021605165
$ws.Range("I25").Value2 = 0.000159715132121236
$ws.Range("J25").Value2 = 0.00004589854229379761
$ws.Range("K25").Value2 = 0.0001028068764918532
$ws.Range("L25").Value2 = 0.0003697555161278859
$ws.Range("M25").Value2 = 0.006346700317619237
$ws.Range("N25").Value2 = 1.000049316702524
$ws.Range("O25").Value2 = 0.006616892489621127
$ws.Range("P25").Value2 = 110.2392809446208
$ws.Range("Q25").Value2 = 165.0886930636898

$ws.Range("A26").Value2 = "model_26_4_0"
$ws.Range("B26").Value2 = 0.9999474700067267
$ws.Range("C26").Value2 = 0.9988418841332297
$ws.Range("D26").Value2 = 0.9991684127059427
$ws.Range("E26").Value2 = 0.9999288340692507
$ws.Range("F26").Value2 = 0.9998797703709493
$ws.Range("G26").Value2 = 0.00004903444292043912
$ws.Range("H26").Value2 = 0.001081050326218786
$ws.Range("I26").Value2 = 0.0001985040908028137
$ws.Range("J26").Value2 = 0.00005130297604605744
$ws.Range("K26").Value2 = 0.0001248969247037952
$ws.Range("L26").Value2 = 0.0003867015158545302
$ws.Range("M26").Value2 = 0.007002459776424219
$ws.Range("N26").Value2 = 1.000060034278027
$ws.Range("O26").Value2 = 0.007300568986196653
$ws.Range("P26").Value2 = 109.8459751800805
$ws.Range("Q26").Value2 = 164.6953872991496
